$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pre Experimental Phase (column C) raw symptom scores, rows 2-23
$values = @(2,0,0,1,1,0,3,0,3,1,2,3,1,0,1,1,0,1,1,4,2,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Cluster / total formulas in column C, mirroring the existing column B formulas
$ws.Range("C24").Formula = "=SUM(C2:C23)"
$ws.Range("C25").Formula = "=SUM(C2 + C12 + C13)"
$ws.Range("C26").Formula = "=SUM(C19:C22)"
$ws.Range("C27").Formula = "=SUM(C14:C18)"
$ws.Range("C28").Formula = "=C23"
$ws.Range("C29").Formula = "=SUM(C3:C6)"
$ws.Range("C30").Formula = "=SUM(C7:C11)"

# Selection moves to C24 as seen in the saved view state
$ws.Range("C24").Select()
